$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# New identifiers / values that replace the old ones throughout the workbook
$oldGuid1 = "4a1395db-5321-4fac-8bad-e0cd38f24991"
$oldGuid2 = "8be85514-958e-44cb-b01f-1ce8cc2cadd6"
$newGuid1 = "df11907f-0aea-4bd0-9a3d-fd2a40b6678b"
$newGuid2 = "ffff24246242-7364-499e-8757-5afc5eb0d618"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"
$newPath1 = "e2e\$newFile1"
$newPath2 = "e2e\$newFile2"

$newXlf1zh = "$newGuid1.0b92c6b94efd9d28c696d9f7b188aa9549738b07.zh-cn.xlf"
$newXlf1de = "$newGuid1.0b92c6b94efd9d28c696d9f7b188aa9549738b07.de-de.xlf"

$newDate1 = "2016-09-02 19:10:44"   # Latest HO Xliff Generate Date / zh H column / de H+K columns (overview/H)
$zhStart  = "2016-09-02 19:10:39"  # zh-cn Correspond Handoff Datetime (H2/H3)
$zhEnd    = "2016-09-02 19:10:57"  # zh-cn Correspond Handback DateTime (K2/K3)
$deStart  = "2016-09-02 19:10:44"  # de-de Correspond Handoff Datetime (H2/H3)
$deEnd    = "2016-09-02 19:11:13"  # de-de Correspond Handback DateTime (K2/K3)

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------

# Plain (non-hyperlinked) cell updates
$ws1.Range("A2").Value = $newFile1
$ws1.Range("A3").Value = $newFile2
$ws1.Range("G2").Value = $newDate1
$ws1.Range("G3").Value = $newDate1

# Hyperlinked cells (B2, B3) - remove + recreate to refresh the display text
# while keeping the same link target.
$hlB2Address = $ws1.Range("B2").Hyperlinks.Item(1).Address
$hlB3Address = $ws1.Range("B3").Hyperlinks.Item(1).Address
if ([string]::IsNullOrEmpty($hlB2Address)) {
    $hlB2Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid1.md"
}
if ([string]::IsNullOrEmpty($hlB3Address)) {
    $hlB3Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid2.md"
}

$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Range("B3").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("B2"), $hlB2Address, "", "", $newPath1)
$ws1.Hyperlinks.Add($ws1.Range("B3"), $hlB3Address, "", "", $newPath2)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------

$ws2.Range("G2").Value = $newXlf1zh
$ws2.Range("H2").Value = $zhStart
$ws2.Range("J2").Value = $newXlf1zh
$ws2.Range("K2").Value = $zhEnd

$ws2.Range("G3").Value = $newXlf1zh
$ws2.Range("H3").Value = $zhStart
$ws2.Range("J3").Value = $newXlf1zh
$ws2.Range("K3").Value = $zhEnd

$hl2A2 = $ws2.Range("A2").Hyperlinks.Item(1).Address
$hl2I2 = $ws2.Range("I2").Hyperlinks.Item(1).Address
$hl2A3 = $ws2.Range("A3").Hyperlinks.Item(1).Address
$hl2I3 = $ws2.Range("I3").Hyperlinks.Item(1).Address
if ([string]::IsNullOrEmpty($hl2A2)) { $hl2A2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid1.md" }
if ([string]::IsNullOrEmpty($hl2I2)) { $hl2I2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4a0639cd7a4cd07b1f8b25dbca2d34667b430229/e2e/$oldGuid1.md" }
if ([string]::IsNullOrEmpty($hl2A3)) { $hl2A3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid2.md" }
if ([string]::IsNullOrEmpty($hl2I3)) { $hl2I3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4a0639cd7a4cd07b1f8b25dbca2d34667b430229/e2e/$oldGuid2.md" }

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Range("I2").Hyperlinks.Delete()
$ws2.Range("A3").Hyperlinks.Delete()
$ws2.Range("I3").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), $hl2A2, "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("I2"), $hl2I2, "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $hl2A3, "", "", $newFile2)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $hl2I3, "", "", $newFile2)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------

$ws3.Range("G2").Value = $newXlf1de
$ws3.Range("H2").Value = $deStart
$ws3.Range("J2").Value = $newXlf1de
$ws3.Range("K2").Value = $deEnd

$ws3.Range("G3").Value = $newXlf1de
$ws3.Range("H3").Value = $deStart
$ws3.Range("J3").Value = $newXlf1de
$ws3.Range("K3").Value = $deEnd

$hl3A2 = $ws3.Range("A2").Hyperlinks.Item(1).Address
$hl3I2 = $ws3.Range("I2").Hyperlinks.Item(1).Address
$hl3A3 = $ws3.Range("A3").Hyperlinks.Item(1).Address
$hl3I3 = $ws3.Range("I3").Hyperlinks.Item(1).Address
if ([string]::IsNullOrEmpty($hl3A2)) { $hl3A2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid1.md" }
if ([string]::IsNullOrEmpty($hl3I2)) { $hl3I2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b7d3a84f0e62556972f1f6156b846de319ad01ca/e2e/$oldGuid1.md" }
if ([string]::IsNullOrEmpty($hl3A3)) { $hl3A3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68a2c554c37d90daae4c9eb058aa0f9292162b12/e2e/$oldGuid2.md" }
if ([string]::IsNullOrEmpty($hl3I3)) { $hl3I3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b7d3a84f0e62556972f1f6156b846de319ad01ca/e2e/$oldGuid2.md" }

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Range("I2").Hyperlinks.Delete()
$ws3.Range("A3").Hyperlinks.Delete()
$ws3.Range("I3").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), $hl3A2, "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("I2"), $hl3I2, "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $hl3A3, "", "", $newFile2)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $hl3I3, "", "", $newFile2)
